# Atualização de bases das ligas, do dia: 19-06-2024 às 21:51
# Swap the full data (columns B..AD) between two pairs of rows.
# Column A (the sequential "id" ordinal) stays put in both pairs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param(
        [int]$Row1,
        [int]$Row2
    )

    $range1 = $ws.Range("B$Row1`:AD$Row1")
    $range2 = $ws.Range("B$Row2`:AD$Row2")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}

Swap-RowData 148 150
Swap-RowData 214 215
